$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price/volume data (and a few name/link swaps reflecting rank changes)

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '31.098.33'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.72%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.957.63'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.76%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.76%  '

# Row 6
$ws.Range('E6').Value = '  +0.12%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4905'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.30%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.71'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.33%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2967'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.21%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06840'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.57%  '

# Row 11
$ws.Range('E11').Value = '  -1.61%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '106.82'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.45%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.943.93'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.14%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07750'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.38%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.420'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.21%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.7127'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.56%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '282.12'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.07%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '31.035.41'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.49%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007757'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.09%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.41%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.199.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.51%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.549'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.48%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.20%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.591'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.98%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.937'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.69%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.50%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.96'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.28%  '

# Row 29
$ws.Range('E29').Value = '  +3.95%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.1057'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.18%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.442'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.15%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.760'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +17.16%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.509'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.80%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05000'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.94%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7673'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.20%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.166'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.22%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02050'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.85%  '

# Row 38
$ws.Range('E38').Value = '  +0.79%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.704'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.53%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.147'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.99%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.446'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.16%  '

# Row 42
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '109.73'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.01%  '

# Row 43
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4477'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.57%  '

# Row 44
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8831'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.67%  '

# Row 45
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '73.03'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.25%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.001'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.03%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.012.12'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +20.54%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.514'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.35%  '

# Row 49
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1270'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.18%  '

# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.400'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.03%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.07'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.15%  '
